# Add two new columns, I ("I0") and J ("IF"), to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1): reuse the same header style already applied to the
# other column headers (s="1": bold, centered, bordered) by copying the
# formatting from the adjacent "IP" header cell, then set the text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J45
$values = @(
    @(4, 4),
    @(9, 9),
    @(5, 6),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(5, 6),
    @(10, 11),
    @(6, 7),
    @(7, 8),
    @(7, 7),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(6, 7),
    @(9, 9),
    @(7, 7),
    @(8, 9),
    @(8, 8),
    @(8, 9),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(6, 7),
    @(10, 10),
    @(7, 7),
    @(7, 7),
    @(6, 6),
    @(5, 5),
    @(7, 7),
    @(4, 4),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(8, 8)
)

for ($idx = 0; $idx -lt $values.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $values[$idx][0]
    $ws.Cells.Item($row, 10).Value = $values[$idx][1]
}
